$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F30").Value = "lin reg"
$ws.Range("F31").Value = "reg gerarchica"
$ws.Range("F32").Value = "reg gerarchica, valutazione didattica, Master DS, possiblità tesi du ANN bayesiane e gerarchiche"
$ws.Range("F33").Value = "revisione progetti"

$ws.Range("F34").Select()
